# Reports upload feature: rename sheet, apply 18pt font styling, and
# populate the word/translation table (rows 2-4) with the new dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name: "Лист1" -> "Лист 1" ---
$ws.Name = "Лист 1"

# --- Apply the body font (18pt, regular) to the whole used area first so
#     the non-bold 18pt font is registered before the bold header font,
#     matching fontId order 1 (plain) / 2 (bold) in the target styles. ---
$ws.Columns("A:D").Font.Size = 18

# --- Header row (row 1): bold 18pt ---
$ws.Range("A1:D1").Font.Bold = $true

# --- Column widths (approximate best-fit widths captured by Excel) ---
$ws.Columns("A").ColumnWidth = 84.83072916666667
$ws.Columns("B").ColumnWidth = 19.385416666666668
$ws.Columns("C").ColumnWidth = 21.721354166666668
$ws.Columns("D").ColumnWidth = 45.830729166666664

# --- Row 2: shind / Англійська / Прикметник / блискучий, сяйливий, лискучий ---
$ws.Range("A2").Value = "shind"
$ws.Range("B2").Value = "Англійська"
$ws.Range("C2").Value = "Прикметник"
$ws.Range("D2").Value = "блискучий, сяйливий, лискучий"

# --- Row 3: dd / Англійська / Прийменник / з, від ---
$ws.Range("A3").Value = "dd"
$ws.Range("B3").Value = "Англійська"
$ws.Range("C3").Value = "Прийменник"
$ws.Range("D3").Value = "з, від"

# --- Row 4 (new row): test / Англійська / Прикметник / ухх, тест, теса ---
$ws.Range("A4").Value = "test"
$ws.Range("B4").Value = "Англійська"
$ws.Range("C4").Value = "Прикметник"
$ws.Range("D4").Value = "ухх, тест, теса"

# --- Make sure the new row also carries the 18pt body font/column width ---
$ws.Range("A4:D4").Font.Size = 18

# --- Selection moves to A4 after the edits, and page orientation is set
#     to portrait (matches the added <pageSetup orientation="portrait".../>) ---
$ws.Range("A4").Select()
$ws.PageSetup.Orientation = 1
